$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 32.935331
$ws.Cells.Item(2, 8).Value = 98.805993
$ws.Cells.Item(2, 9).Value = 0.1836164637112342
$ws.Cells.Item(2, 10).Value = 0.1836164637112342
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 0.7971666666666667
$ws.Cells.Item(2, 14).Value = 2.3915
$ws.Cells.Item(2, 15).Value = 0.3308932235309289
$ws.Cells.Item(2, 16).Value = 0.3308932235309289
$ws.Cells.Item(2, 17).Value = 26.25494802883333
$ws.Cells.Item(2, 18).Value = 236.2945322595
$ws.Cells.Item(2, 19).Value = 0.06075744357076009
$ws.Cells.Item(2, 20).Value = 0.06075744357076009

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 32.935331
$ws.Cells.Item(3, 8).Value = 98.805993
$ws.Cells.Item(3, 9).Value = 0.1836164637112342
$ws.Cells.Item(3, 10).Value = 0.1836164637112342
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 0.679891
$ws.Cells.Item(3, 14).Value = 2.039673
$ws.Cells.Item(3, 15).Value = 0.282213662521012
$ws.Cells.Item(3, 16).Value = 0.282213662521012
$ws.Cells.Item(3, 17).Value = 22.392435128921
$ws.Cells.Item(3, 18).Value = 201.531916160289
$ws.Cells.Item(3, 19).Value = 0.05181907472310389
$ws.Cells.Item(3, 20).Value = 0.05181907472310389

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 32.935331
$ws.Cells.Item(4, 8).Value = 98.805993
$ws.Cells.Item(4, 9).Value = 0.1836164637112342
$ws.Cells.Item(4, 10).Value = 0.1836164637112342
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.76147
$ws.Cells.Item(4, 14).Value = 2.28441
$ws.Cells.Item(4, 15).Value = 0.3160760145374406
$ws.Cells.Item(4, 16).Value = 0.3160760145374406
$ws.Cells.Item(4, 17).Value = 25.07926649657
$ws.Cells.Item(4, 18).Value = 225.71339846913
$ws.Cells.Item(4, 19).Value = 0.05803676005330548
$ws.Cells.Item(4, 20).Value = 0.05803676005330548

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 32.935331
$ws.Cells.Item(5, 8).Value = 98.805993
$ws.Cells.Item(5, 9).Value = 0.1836164637112342
$ws.Cells.Item(5, 10).Value = 0.1836164637112342
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 0.170608
$ws.Cells.Item(5, 14).Value = 0.5118240000000001
$ws.Cells.Item(5, 15).Value = 0.07081709941061849
$ws.Cells.Item(5, 16).Value = 0.07081709941061851
$ws.Cells.Item(5, 17).Value = 5.619030951248
$ws.Cells.Item(5, 18).Value = 50.57127856123201
$ws.Cells.Item(5, 19).Value = 0.01300318536406469
$ws.Cells.Item(5, 20).Value = 0.0130031853640647

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 95.562134
$ws.Cells.Item(6, 8).Value = 286.686402
$ws.Cells.Item(6, 9).Value = 0.5327646808765668
$ws.Cells.Item(6, 10).Value = 0.5327646808765667
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 0.7971666666666667
$ws.Cells.Item(6, 14).Value = 2.3915
$ws.Cells.Item(6, 15).Value = 0.3308932235309289
$ws.Cells.Item(6, 16).Value = 0.3308932235309289
$ws.Cells.Item(6, 17).Value = 76.17894782033333
$ws.Cells.Item(6, 18).Value = 685.610530383
$ws.Cells.Item(6, 19).Value = 0.1762882226386738
$ws.Cells.Item(6, 20).Value = 0.1762882226386738

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 95.562134
$ws.Cells.Item(7, 8).Value = 286.686402
$ws.Cells.Item(7, 9).Value = 0.5327646808765668
$ws.Cells.Item(7, 10).Value = 0.5327646808765667
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 0.679891
$ws.Cells.Item(7, 14).Value = 2.039673
$ws.Cells.Item(7, 15).Value = 0.282213662521012
$ws.Cells.Item(7, 16).Value = 0.282213662521012
$ws.Cells.Item(7, 17).Value = 64.97183484739401
$ws.Cells.Item(7, 18).Value = 584.7465136265459
$ws.Cells.Item(7, 19).Value = 0.1503534718520141
$ws.Cells.Item(7, 20).Value = 0.1503534718520141

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 95.562134
$ws.Cells.Item(8, 8).Value = 286.686402
$ws.Cells.Item(8, 9).Value = 0.5327646808765668
$ws.Cells.Item(8, 10).Value = 0.5327646808765667
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 0.76147
$ws.Cells.Item(8, 14).Value = 2.28441
$ws.Cells.Item(8, 15).Value = 0.3160760145374406
$ws.Cells.Item(8, 16).Value = 0.3160760145374406
$ws.Cells.Item(8, 17).Value = 72.76769817698001
$ws.Cells.Item(8, 18).Value = 654.90928359282
$ws.Cells.Item(8, 19).Value = 0.1683941370177766
$ws.Cells.Item(8, 20).Value = 0.1683941370177766

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 95.562134
$ws.Cells.Item(9, 8).Value = 286.686402
$ws.Cells.Item(9, 9).Value = 0.5327646808765668
$ws.Cells.Item(9, 10).Value = 0.5327646808765667
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.170608
$ws.Cells.Item(9, 14).Value = 0.5118240000000001
$ws.Cells.Item(9, 15).Value = 0.07081709941061849
$ws.Cells.Item(9, 16).Value = 0.07081709941061851
$ws.Cells.Item(9, 17).Value = 16.303664557472
$ws.Cells.Item(9, 18).Value = 146.732981017248
$ws.Cells.Item(9, 19).Value = 0.03772884936810227
$ws.Cells.Item(9, 20).Value = 0.03772884936810227

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 37.06916133333333
$ws.Cells.Item(10, 8).Value = 111.207484
$ws.Cells.Item(10, 9).Value = 0.2066628180165514
$ws.Cells.Item(10, 10).Value = 0.2066628180165514
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 0.7971666666666667
$ws.Cells.Item(10, 14).Value = 2.3915
$ws.Cells.Item(10, 15).Value = 0.3308932235309289
$ws.Cells.Item(10, 16).Value = 0.3308932235309289
$ws.Cells.Item(10, 17).Value = 29.55029977622223
$ws.Cells.Item(10, 18).Value = 265.952697986
$ws.Cells.Item(10, 19).Value = 0.06838332603748243
$ws.Cells.Item(10, 20).Value = 0.06838332603748241

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 37.06916133333333
$ws.Cells.Item(11, 8).Value = 111.207484
$ws.Cells.Item(11, 9).Value = 0.2066628180165514
$ws.Cells.Item(11, 10).Value = 0.2066628180165514
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 0.679891
$ws.Cells.Item(11, 14).Value = 2.039673
$ws.Cells.Item(11, 15).Value = 0.282213662521012
$ws.Cells.Item(11, 16).Value = 0.282213662521012
$ws.Cells.Item(11, 17).Value = 25.20298916808133
$ws.Cells.Item(11, 18).Value = 226.826902512732
$ws.Cells.Item(11, 19).Value = 0.05832307077936437
$ws.Cells.Item(11, 20).Value = 0.05832307077936436

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 37.06916133333333
$ws.Cells.Item(12, 8).Value = 111.207484
$ws.Cells.Item(12, 9).Value = 0.2066628180165514
$ws.Cells.Item(12, 10).Value = 0.2066628180165514
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 0.76147
$ws.Cells.Item(12, 14).Value = 2.28441
$ws.Cells.Item(12, 15).Value = 0.3160760145374406
$ws.Cells.Item(12, 16).Value = 0.3160760145374406
$ws.Cells.Item(12, 17).Value = 28.22705428049333
$ws.Cells.Item(12, 18).Value = 254.04348852444
$ws.Cells.Item(12, 19).Value = 0.06532115987174796
$ws.Cells.Item(12, 20).Value = 0.06532115987174793

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 37.06916133333333
$ws.Cells.Item(13, 8).Value = 111.207484
$ws.Cells.Item(13, 9).Value = 0.2066628180165514
$ws.Cells.Item(13, 10).Value = 0.2066628180165514
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 0.170608
$ws.Cells.Item(13, 14).Value = 0.5118240000000001
$ws.Cells.Item(13, 15).Value = 0.07081709941061849
$ws.Cells.Item(13, 16).Value = 0.07081709941061851
$ws.Cells.Item(13, 17).Value = 6.324295476757333
$ws.Cells.Item(13, 18).Value = 56.918659290816
$ws.Cells.Item(13, 19).Value = 0.01463526132795668
$ws.Cells.Item(13, 20).Value = 0.01463526132795668

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 13.80362366666667
$ws.Cells.Item(14, 8).Value = 41.410871
$ws.Cells.Item(14, 9).Value = 0.07695603739564764
$ws.Cells.Item(14, 10).Value = 0.07695603739564763
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 0.7971666666666667
$ws.Cells.Item(14, 14).Value = 2.3915
$ws.Cells.Item(14, 15).Value = 0.3308932235309289
$ws.Cells.Item(14, 16).Value = 0.3308932235309289
$ws.Cells.Item(14, 17).Value = 11.00378866627778
$ws.Cells.Item(14, 18).Value = 99.03409799650001
$ws.Cells.Item(14, 19).Value = 0.02546423128401256
$ws.Cells.Item(14, 20).Value = 0.02546423128401255

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 13.80362366666667
$ws.Cells.Item(15, 8).Value = 41.410871
$ws.Cells.Item(15, 9).Value = 0.07695603739564764
$ws.Cells.Item(15, 10).Value = 0.07695603739564763
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 0.679891
$ws.Cells.Item(15, 14).Value = 2.039673
$ws.Cells.Item(15, 15).Value = 0.282213662521012
$ws.Cells.Item(15, 16).Value = 0.282213662521012
$ws.Cells.Item(15, 17).Value = 9.384959498353666
$ws.Cells.Item(15, 18).Value = 84.46463548518301
$ws.Cells.Item(15, 19).Value = 0.02171804516652969
$ws.Cells.Item(15, 20).Value = 0.02171804516652968

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 13.80362366666667
$ws.Cells.Item(16, 8).Value = 41.410871
$ws.Cells.Item(16, 9).Value = 0.07695603739564764
$ws.Cells.Item(16, 10).Value = 0.07695603739564763
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 0.76147
$ws.Cells.Item(16, 14).Value = 2.28441
$ws.Cells.Item(16, 15).Value = 0.3160760145374406
$ws.Cells.Item(16, 16).Value = 0.3160760145374406
$ws.Cells.Item(16, 17).Value = 10.51104531345667
$ws.Cells.Item(16, 18).Value = 94.59940782110999
$ws.Cells.Item(16, 19).Value = 0.02432395759461055
$ws.Cells.Item(16, 20).Value = 0.02432395759461054

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 13.80362366666667
$ws.Cells.Item(17, 8).Value = 41.410871
$ws.Cells.Item(17, 9).Value = 0.07695603739564764
$ws.Cells.Item(17, 10).Value = 0.07695603739564763
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 0.170608
$ws.Cells.Item(17, 14).Value = 0.5118240000000001
$ws.Cells.Item(17, 15).Value = 0.07081709941061849
$ws.Cells.Item(17, 16).Value = 0.07081709941061851
$ws.Cells.Item(17, 17).Value = 2.355008626522667
$ws.Cells.Item(17, 18).Value = 21.195077638704
$ws.Cells.Item(17, 19).Value = 0.005449803350494853
$ws.Cells.Item(17, 20).Value = 0.005449803350494853
